# Separated languages to make it easier to add more languages later
#
# Adds a French ("fr") column/translation alongside the existing English
# content:
#   - Sheet1:            new row A3 = "en,fr" (list of languages present)
#   - cutscene01 (Sh2):  new column C with "fr" header + translated rows
#   - cutscene02..06:    new column C with just the "fr" header (no
#                         translations added yet for those cutscenes)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# cutscene01 (sheet index 2): full French translation column added.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C1").Value = "fr"
$ws2.Range("C1").Font.Bold = $true

$ws2.Range("C2").Value = "Lorsque les moines ont commencé à mourir d'une maladie mystérieuse, le Dr Ulko est arrivé sur les lieux."
$ws2.Range("C2").WrapText = $true

$ws2.Range("C3").Value = "Malheureusement pour lui, les moines ont pris l'histoire à l'envers et ont cru que la présence de l'étranger avait causé la maladie, et non que la maladie avait causé la présence de l'étranger. Personne ne sait s'ils avaient raison ou non."
$ws2.Range("C3").WrapText = $true

$ws2.Range("C4").Value = "Ils décident de l'exiler dans l'Autre Monde..."
$ws2.Range("C4").WrapText = $true

$ws2.Columns.Item(3).ColumnWidth = 58.25

$ws2.Rows.Item(1).RowHeight = 18.75
$ws2.Rows.Item(2).RowHeight = 49.5

# ---------------------------------------------------------------------
# Sheet1 (sheet index 1): records which languages are present.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A3").Value = "en,fr"

# ---------------------------------------------------------------------
# cutscene02..cutscene06 (sheet indices 3-7): just add the "fr" header
# cell in column C so the sheet is ready for future translations.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C1").Value = "fr"
$ws3.Range("C1").Font.Bold = $true

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("C1").Value = "fr"
$ws4.Range("C1").Font.Bold = $true

$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("C1").Value = "fr"
$ws5.Range("C1").Font.Bold = $true

$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("C1").Value = "fr"
$ws6.Range("C1").Font.Bold = $true

$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("C1").Value = "fr"
$ws7.Range("C1").Font.Bold = $true

# ---------------------------------------------------------------------
# Selections, matching where the author's cursor ended up on each sheet.
# Applied in tab order, cutscene06 last, so the active tab stays on
# cutscene06 (index 6) like in the original workbook.
# ---------------------------------------------------------------------
$ws1.Range("A3").Select() | Out-Null
$ws3.Range("C1").Select() | Out-Null
$ws4.Range("C1").Select() | Out-Null
$ws5.Range("C1").Select() | Out-Null
$ws6.Range("C1").Select() | Out-Null
$ws7.Range("D10").Select() | Out-Null
